$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.979.95"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "2.269.35"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "317.98"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "102.74"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "38.64"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").Value = "0.0839"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "7.84"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "2.619.23"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "0.876"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "2.286.55"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "43.943.20"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "14.50"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "6.65"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "66.08"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "238.66"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "4.05"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "39.29"
$ws.Range("E29").Value = "  +15.14%  "
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "6.56"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "162.42"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0882"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "20.46"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "2.72"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "3.32"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "2.01"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "4.56"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "0.107"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  +6.93%  "
$ws.Range("D42").Value = "16.14"
$ws.Range("E42").Value = "  +31.84%  "
$ws.Range("D43").Value = "0.0326"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "1.789.05"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").Value = "0.208"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "85.41"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").Value = "5.40"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  +5.60%  "
$ws.Range("D50").Value = "75.79"
$ws.Range("E50").Value = "  -5.35%  "
$ws.Range("D51").Value = "59.72"
$ws.Range("E51").Value = "  -2.93%  "
